{"js": "// Apply the ANTS Methodology Summary edits described by the commit\n// \"adding comments on gray matter probability\".\n//\n// Strategy: use Range.search() to locate unique text anchors and\n// Range.insertText(..., \"Replace\") to rewrite them in place, plus a\n// paragraph-level firstLineIndent tweak and a trailing blank-paragraph\n// removal.\n\n// 1) Give the title paragraph a first-line indent of 12pt\n//    (OOXML <w:ind w:firstLine=\"240\"/>, 240 twips == 12 pt).\nconst titlePara = context.document.body.paragraphs.getFirst();\ntitlePara.firstLineIndent = 12;\n\n// Helper: replace the first (expected unique) occurrence of `findText`\n// with `replaceText`, preserving the formatting of the text it replaces.\nasync function replaceOnce(findText, replaceText) {\n  const results = context.document.body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly one match for \" + JSON.stringify(findText) +\n      \" but found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"The core processing is based on mapping T1\" -> \"The core processing maps T1\"\nawait replaceOnce(\n  \"processing is based on mapping T1\",\n  \"processing maps T1\"\n);\n\n// 3) \"maps are able to capture\" -> \"maps capture\"\nawait replaceOnce(\n  \" are able to capture\",\n  \" capture\"\n);\n\n// 4) \"is perhaps the most stable\" -> \"is the most stable\"\nawait replaceOnce(\n  \"is perhaps the most stable\",\n  \"is the most stable\"\n);\n\n// 5) \"These algorithms allow template-based priors\" ->\n//    \"After defining the template image to target image coordinate\n//     transformation, we employ template-based priors\"\nawait replaceOnce(\n  \"These algorithms allow template-based priors\",\n  \"After defining the template image to target image coordinate transformation, we employ template-based priors\"\n);\n\n// 6) \"to guide cortical segmentation and compute\" ->\n//    \"to guide cortical gray matter segmentation, cortical parcellation and compute\"\nawait replaceOnce(\n  \"to guide cortical segmentation and compute\",\n  \"to guide cortical gray matter segmentation, cortical parcellation and compute\"\n);\n\n// 7) \"either exact or partial matching [8].\" ->\n//    \"either exact or partial matching between labels or landmarks defined by the user [8].\"\nawait replaceOnce(\n  \"either exact or partial matching [8].\",\n  \"either exact or partial matching between labels or landmarks defined by the user [8].\"\n);\n\n// 8) \"ANTS may be used to consistently standardize and normalize data from T1\" ->\n//    \"ANTS is applicable to T1\"\nawait replaceOnce(\n  \"ANTS may be used to consistently standardize and normalize data from T1\",\n  \"ANTS is applicable to T1\"\n);\n\n// 9) \"T1.  The use of both modalities together optimizes the normalization\" ->\n//    \"T1. We use both modalities together to optimize the normalization\"\nawait replaceOnce(\n  \"T1.  The use of both modalities together optimizes the normalization\",\n  \"T1. We use both modalities together to optimize the normalization\"\n);\n\n// 10) \"to reconstruct both cortical thickness\" -> \"to reconstruct gray matter, cortical thickness\"\nawait replaceOnce(\n  \"to reconstruct both cortical thickness\",\n  \"to reconstruct gray matter, cortical thickness\"\n);\n\n// 11) \"both traditional jacobian studies\" -> \"both traditional Jacobian-based morphometry studies\"\nawait replaceOnce(\n  \"both traditional jacobian studies\",\n  \"both traditional Jacobian-based morphometry studies\"\n);\n\n// 12) \"optimal power for both gray (cortical thickness)\" ->\n//     \"optimal power for both cortex (gray matter probability, cortical thickness)\"\nawait replaceOnce(\n  \"optimal power for both gray (cortical thickness)\",\n  \"optimal power for both cortex (gray matter probability, cortical thickness)\"\n);\n\n// 13) Remove the extra blank paragraph that used to sit between the\n//     [7] Avants reference and the final blank paragraph before sectPr.\nconst refResults = context.document.body.search(\n  \"Landmark and Intensity Driven Lagrangian\",\n  { matchCase: true }\n);\nrefResults.load(\"items\");\nawait context.sync();\nconst refPara = refResults.items[0].paragraphs.getFirst();\nconst blankPara = refPara.getNext();\nblankPara.delete();\n\nawait context.sync();\n", "ps1": "# Apply the ANTS Methodology Summary edits described by the commit\n# \"adding comments on gray matter probability\".\n#\n# Strategy: use Range.Find/Replace (Execute with Replace:=wdReplaceOne,\n# value 1) to rewrite unique text anchors in place, plus a paragraph-level\n# FirstLineIndent tweak and a trailing blank-paragraph removal.\n\n$d = $word.ActiveDocument\n\n# 1) Give the title paragraph a first-line indent of 12pt\n#    (OOXML <w:ind w:firstLine=\"240\"/>, 240 twips == 12 pt).\n$d.Paragraphs(1).Range.ParagraphFormat.FirstLineIndent = 12\n\nfunction ReplaceOnce($findText, $replaceText) {\n    $r = $d.Content\n    $f = $r.Find\n    $f.ClearFormatting()\n    $f.Replacement.ClearFormatting()\n    $result = $f.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $result) {\n        throw \"Replacement not found: $findText\"\n    }\n}\n\n# 2) \"The core processing is based on mapping T1\" -> \"The core processing maps T1\"\nReplaceOnce \"processing is based on mapping T1\" \"processing maps T1\"\n\n# 3) \"maps are able to capture\" -> \"maps capture\"\nReplaceOnce \" are able to capture\" \" capture\"\n\n# 4) \"is perhaps the most stable\" -> \"is the most stable\"\nReplaceOnce \"is perhaps the most stable\" \"is the most stable\"\n\n# 5) \"These algorithms allow template-based priors\" ->\n#    \"After defining the template image to target image coordinate\n#     transformation, we employ template-based priors\"\nReplaceOnce \"These algorithms allow template-based priors\" \"After defining the template image to target image coordinate transformation, we employ template-based priors\"\n\n# 6) \"to guide cortical segmentation and compute\" ->\n#    \"to guide cortical gray matter segmentation, cortical parcellation and compute\"\nReplaceOnce \"to guide cortical segmentation and compute\" \"to guide cortical gray matter segmentation, cortical parcellation and compute\"\n\n# 7) \"either exact or partial matching [8].\" ->\n#    \"either exact or partial matching between labels or landmarks defined by the user [8].\"\nReplaceOnce \"either exact or partial matching [8].\" \"either exact or partial matching between labels or landmarks defined by the user [8].\"\n\n# 8) \"ANTS may be used to consistently standardize and normalize data from T1\" ->\n#    \"ANTS is applicable to T1\"\nReplaceOnce \"ANTS may be used to consistently standardize and normalize data from T1\" \"ANTS is applicable to T1\"\n\n# 9) \"T1.  The use of both modalities together optimizes the normalization\" ->\n#    \"T1. We use both modalities together to optimize the normalization\"\nReplaceOnce \"T1.  The use of both modalities together optimizes the normalization\" \"T1. We use both modalities together to optimize the normalization\"\n\n# 10) \"to reconstruct both cortical thickness\" -> \"to reconstruct gray matter, cortical thickness\"\nReplaceOnce \"to reconstruct both cortical thickness\" \"to reconstruct gray matter, cortical thickness\"\n\n# 11) \"both traditional jacobian studies\" -> \"both traditional Jacobian-based morphometry studies\"\nReplaceOnce \"both traditional jacobian studies\" \"both traditional Jacobian-based morphometry studies\"\n\n# 12) \"optimal power for both gray (cortical thickness)\" ->\n#     \"optimal power for both cortex (gray matter probability, cortical thickness)\"\nReplaceOnce \"optimal power for both gray (cortical thickness)\" \"optimal power for both cortex (gray matter probability, cortical thickness)\"\n\n# 13) Remove the extra blank paragraph that used to sit between the\n#     [7] Avants reference and the final blank paragraph before sectPr.\n$r = $d.Content\n$r.Find.Execute(\"Landmark and Intensity Driven Lagrangian\") | Out-Null\n$refPara = $r.Paragraphs(1)\n$blankPara = $refPara.Next()\n$blankPara.Range.Delete()\n"}
